# Marksheet result fix: recompute Right/Wrong/Not-Attempt/Max and the
# per-question Student Ans columns now that the student's answers are
# known (previously the sheet was an unfilled "Absent" template).
# Also drops the unused third Student/Correct-Ans column pair (G:H) and
# the now-unused D:E rows beyond the first 3 answered questions, since
# only one fill-in-the-blank answer pair has real content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Summary rows (10-12): give the row-label cells in column A the
# same "mtitleStyle" formatting already used by the header row (row 9) ----
$labelCells = @("A10", "A11", "A12")
foreach ($cell in $labelCells) {
    $ws.Range("A9").Copy()
    $ws.Range($cell).PasteSpecial(-4122)
}

# ---- Score summary values ----
$ws.Range("B10").Value = 17
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 8
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 68
$ws.Range("C12").Value = -3
$ws.Range("E12").Value = "65/112"

# ---- Remove the unused third Student Ans / Correct Ans column pair ----
$ws.Range("G15:H21").Clear()

# ---- Remove the now-unused D:E answer rows (only rows 16-18 keep data) ----
$ws.Range("D19:E40").Clear()

# ---- Fill in column D (2nd Student Ans) for rows 16-18, all correct ----
$correctD = @(
    @("D16", "Option A"),
    @("D17", "Option C"),
    @("D18", "Option D")
)
foreach ($pair in $correctD) {
    $ws.Range("B10").Copy()
    $ws.Range($pair[0]).PasteSpecial(-4122)
    $ws.Range($pair[0]).Value = $pair[1]
}

# ---- Fill in column A (1st Student Ans) answers, rows 16-40 ----
# Correctly-answered questions -> "correctStyle" (same formatting as B10)
$correctA = @(
    @("A16", "Option A"),
    @("A18", "Option B"),
    @("A19", "Option C"),
    @("A22", "Option D"),
    @("A25", "Option A"),
    @("A27", "Option A"),
    @("A28", "Option D"),
    @("A29", "Option D"),
    @("A30", "Option B"),
    @("A32", "Option C"),
    @("A33", "Option D"),
    @("A35", "Option D"),
    @("A36", "Option A"),
    @("A39", "Option D")
)
foreach ($pair in $correctA) {
    $ws.Range("B10").Copy()
    $ws.Range($pair[0]).PasteSpecial(-4122)
    $ws.Range($pair[0]).Value = $pair[1]
}

# Incorrectly-answered questions -> "incorrectStyle" (same formatting as C10)
$incorrectA = @(
    @("A20", "Option A"),
    @("A26", "Option D"),
    @("A34", "Option A")
)
foreach ($pair in $incorrectA) {
    $ws.Range("C10").Copy()
    $ws.Range($pair[0]).PasteSpecial(-4122)
    $ws.Range($pair[0]).Value = $pair[1]
}

# Rows 17, 21, 23, 24, 31, 37, 38, 40 in column A remain un-attempted
# (already blank with the "normalStyle" formatting), so nothing to do there.
